# Refresh the "cryptos" price list (Coin / Link / Price / Volume(1h))
# to the latest snapshot pulled by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '41.793.13'
$ws.Cells.Item(2, 5).Value = '  +0.00%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.172.30'
$ws.Cells.Item(3, 5).Value = '  -1.85%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.04%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.Value2 = '''239.20'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.79%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -2.55%  '

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.Value2 = '''71.78'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -1.23%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.04%  '

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.Value2 = '''0.578'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -3.23%  '

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.Value2 = '''40.17'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.76%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.Value2 = '''0.0909'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -3.68%  '

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.Value2 = '''54.40'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -4.05%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.Value2 = '''6.72'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -3.49%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -2.73%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.499.35'
$ws.Cells.Item(15, 5).Value = '  -1.89%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.Value2 = '''14.31'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.93%  '

# Coinranking's page reordered #16/#17 (WrappedEther now ranks
# above Polygon) -- the two rows swap places/content.
# Row 17
$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '2.166.88'
$ws.Cells.Item(17, 5).Value = '  -2.03%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'Polygon'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Cells.Item(18, 4)
$cell.Value2 = '''0.788'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -5.22%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '41.609.96'
$ws.Cells.Item(19, 5).Value = '  -0.21%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -2.82%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.Value2 = '''69.79'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -3.45%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.Value2 = '''5.77'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -6.02%  '

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.Value2 = '''10.09'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -8.49%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.Value2 = '''226.86'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.86%  '

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.Value2 = '''1.98'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -3.12%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.20%  '

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.Value2 = '''10.71'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -5.39%  '

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.Value2 = '''3.32'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -8.36%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.Value2 = '''2.19'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -3.34%  '

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.Value2 = '''171.97'
$cell.Style = "Normal"

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.Value2 = '''19.81'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -2.72%  '

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.Value2 = '''32.52'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +7.94%  '

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.Value2 = '''0.0773'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -2.75%  '

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.Value2 = '''5.13'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -7.81%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -2.86%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -1.22%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +1.79%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.Value2 = '''0.0302'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +1.19%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.Value2 = '''12.55'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -6.86%  '

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.Value2 = '''2.06'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.30%  '

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.Value2 = '''5.35'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -4.54%  '

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.Value2 = '''59.34'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -6.94%  '

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.Value2 = '''0.190'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -3.03%  '

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.Value2 = '''8.32'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -3.81%  '

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.Value2 = '''0.0971'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -2.58%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.Value2 = '''97.82'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -4.73%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.Value2 = '''1.08'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.74%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -4.03%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -6.30%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -2.42%  '
